$d = $word.ActiveDocument

$replacements = @(
    @("2024-09-03 Tuesday", "2024-09-04 Wednesday"),
    @("60×24=", "41×68="),
    @("35×73=", "24×46="),
    @("33×90=", "32×56="),
    @("26×85=", "63×65="),
    @("42×77=", "41×74="),
    @("60×14=", "31×15="),
    @("96×63=", "45×70="),
    @("66×53=", "87×83="),
    @("70×37=", "50×22="),
    @("20×33=", "76×89="),
    @("72×70=", "47×30="),
    @("17×63=", "22×32="),
    @("32×15=", "88×21="),
    @("91×65=", "64×26="),
    @("24×95=", "19×21="),
    @("98×67=", "59×56="),
    @("70×30=", "97×32="),
    @("48×79=", "51×93="),
    @("62×95=", "18×16="),
    @("11×14=", "62×91="),
    @("78×20=", "67×67="),
    @("49×30=", "87×56="),
    @("11×81=", "67×81="),
    @("55×54=", "15×49="),
    @("87×58=", "73×24=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $found = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $found) {
        Write-Host "WARNING: did not find" $old
    }
}

Write-Host "Done applying" $replacements.Count "replacements"
